$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the J-column formulas for rows 9-12
$ws.Range("J9").Formula = "=60+(H9-F9)*E9"
$ws.Range("J10").Formula = "=60+(H10-F10)*E10"
$ws.Range("J11").Formula = "=60+(H11-F11)*E11"
$ws.Range("J12").Formula = "=60+(H12-F12)*E12"

# Column J should now share column I's width/bestFit/customWidth settings
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Sheet view changes: scroll so A4 is the top-left cell, and change selection
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("L11").Select()
